$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '24.398.12'
$ws.Range("E2").Value = '  -2.35%  '

# Row 3
$ws.Range("D3").Value = '1.647.94'
$ws.Range("E3").Value = '  -4.18%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.38%  '

# Row 6
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.16%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3613'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.40%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.65'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.61%  '

# Row 9
$ws.Range("E9").Value = '  -6.32%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.118'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.30%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06914'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.41%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9996'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.05%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.918'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.41%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.19%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.568'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.09%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.647.12'
$ws.Range("E16").Value = '  -4.28%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001038'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -7.80%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06495'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.37%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9997'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.03%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '76.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -9.46%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.903'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.59%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.63'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -9.56%  '

# Row 23
$ws.Range("E23").Value = '  -8.43%  '

# Row 24
$ws.Range("D24").Value = '24.350.67'
$ws.Range("E24").Value = '  -2.52%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.422'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.64%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.320'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -17.25%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '146.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.07%  '

# Row 28
$ws.Range("E28").Value = '  -11.75%  '

# Row 29
$ws.Range("D29").Value = '1.829.59'
$ws.Range("E29").Value = '  -4.28%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.20'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.12%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.163'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.19%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.042'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.27%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.583'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -18.36%  '

# Row 34
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08325'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.70%  '

# Row 35
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.673'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.68%  '

# Row 36
$ws.Range("E36").Value = '  -11.99%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.200'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.72%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06004'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.47%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02194'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.12%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.204'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.13%  '

# Row 41
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2040'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.22%  '

# Row 42
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.149'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.49%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.08%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5783'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -10.44%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.723'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.95%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.05%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5520'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.13%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '121.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.57%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.929'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -10.28%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06880'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.63%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.60%  '
